$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header row: "_old" suffix -> "_FV2310", "_new" suffix -> "_FV2404"
$ws.Range("A1").Value = "Segmentname_FV2310"
$ws.Range("B1").Value = "Segmentgruppe_FV2310"
$ws.Range("C1").Value = "Segment_FV2310"
$ws.Range("D1").Value = "Datenelement_FV2310"
$ws.Range("E1").Value = "Segment ID_FV2310"
$ws.Range("F1").Value = "Code_FV2310"
$ws.Range("G1").Value = "Qualifier_FV2310"
$ws.Range("H1").Value = "Beschreibung_FV2310"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("J1").Value = "Bedingung_FV2310"
$ws.Range("K1").Value = "diff"
$ws.Range("L1").Value = "Segmentname_FV2404"
$ws.Range("M1").Value = "Segmentgruppe_FV2404"
$ws.Range("N1").Value = "Segment_FV2404"
$ws.Range("O1").Value = "Datenelement_FV2404"
$ws.Range("P1").Value = "Segment ID_FV2404"
$ws.Range("Q1").Value = "Code_FV2404"
$ws.Range("R1").Value = "Qualifier_FV2404"
$ws.Range("S1").Value = "Beschreibung_FV2404"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("U1").Value = "Bedingung_FV2404"

# 2) Turn the used range into a proper Excel table (adds xl/tables/table1.xml
#    and the <tableParts> reference on the worksheet).
#    The engine bakes the header row's *current* cell formatting into a new
#    dxf (headerRowDxfId) whenever ListObjects.Add sees non-default styling
#    on row 1. Row 1 already carries the workbook's bold/centered header
#    style (style index 1), so: stash that formatting on a scratch cell,
#    blank the header's style out, create the table (now against a
#    "default-styled" header so no dxf is recorded), then paste the
#    original formatting back and drop the scratch row again.
$hdr = $ws.Range("A1:U1")
$scratch = $ws.Range("A200")
$hdr.Copy()
$scratch.PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$hdr.Style = "Normal"

$rng = $ws.Range("A1:U79")
$lo = $ws.ListObjects.Add(1, $rng, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

$scratch.Copy()
$hdr.PasteSpecial(-4122)       # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Rows(200).Delete()

# 3) Freeze the header row (split below row 1, frozen).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
